# Planification.xlsx update
# - Fill in actuals for "Recherches liées au projet" (row 6, S0 & s1)
# - Clear the stray actual on "Rechercher le matériel sur internet" (row 19, S0)
# - Update actuals for "Rédaction du cahier des charges" (row 39, s1) and add s2
# - Add actuals for "Séances ébdomadaires avec prof. Répondant" (row 41, s2)
# - Add actuals for "Séances ébdomadaires avec prof. Répondant (préparation)" (row 42, s2)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 6 - T10.40 Recherches liées au projet
$ws.Range("E6").Value = 1.5
$ws.Range("F6").Value = 5

# Row 19 - T40.20 Rechercher le matériel sur internet
$ws.Range("E19").Value = ""

# Row 39 - T80.20 Rédaction du cahier des charges
$ws.Range("F39").Value = 6
$ws.Range("G39").Value = 3

# Row 41 - T80.40 Séances ébdomadaires avec prof. Répondant
$ws.Range("G41").Value = 1

# Row 42 - T80.50 Séances ébdomadaires avec prof. Répondant (préparation)
$ws.Range("G42").Value = 0.25

# Update the view: zoom level and current selection
$ws.Activate()
$excel.ActiveWindow.Zoom = 58
$ws.Range("G40").Select()
